# Driver file - result sheet is added with one more column failed functions
#
# Login sheet:
#   - A3 ("Login ID" row) now holds a new test login email, with the
#     previous email moved down into a brand new row 12 (mirroring the old
#     hyperlink/style).
#   - Login becomes the active sheet/tab with A3 selected.

$wb = $excel.ActiveWorkbook
$wsLogin  = $wb.Worksheets.Item("Login")

# A3 used to read "test16@styletag.com" -- replace it with the new address.
$wsLogin.Range("A3").Value = "test3456@styletag.com"

# Re-home the old address on a new row 12, keeping the same look (style)
# as A3 had, then wire up its hyperlink.
$wsLogin.Range("A12").Style = $wsLogin.Range("A3").Style
$wsLogin.Hyperlinks.Add(
    $wsLogin.Range("A12"),
    "mailto:test16@styletag.com",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "test16@styletag.com"
)

# Login is now the selected/active sheet (was SearchKey before), with A3
# the active cell (was D11 before).
$null = $wsLogin.Select()
$null = $wsLogin.Range("A3").Select()
